$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to Text format first, so numeric-looking
# strings (e.g. "5.176", "9.198", "1.002") are preserved exactly as
# text and are not auto-converted into numbers by Excel.
$cells = @("D2","D3","E3","E4","D5","E5","D6","E6","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","B13","C13","D13","E13","B14","C14","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","E23","D24","E24","D25","E25","D26","E26","D27","E27","E28","D29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","E35","D36","E36","D37","E37","D38","E38","D39","E39","E40","D41","E41","B42","C42","D42","E42","B43","C43","D43","E43","D44","E44","D45","E45","E46","D47","E47","B48","C48","D48","E48","B49","C49","D49","E49","B50","C50","D50","E50","B51","C51","D51","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "29.089.33"
$ws.Range("D3").Value = "1.842.00"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "0.6909"
$ws.Range("E5").Value = "  -6.81%  "
$ws.Range("D6").Value = "236.60"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.3028"
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").Value = "0.07506"
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("D10").Value = "23.24"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("D11").Value = "0.08078"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").Value = "0.7196"
$ws.Range("E12").Value = "  -4.05%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.822.75"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.176"
$ws.Range("E14").Value = "  -3.86%  "
$ws.Range("D15").Value = "88.63"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").Value = "29.320.75"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "5.769"
$ws.Range("E17").Value = "  -5.67%  "
$ws.Range("D18").Value = "240.52"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").Value = "0.000007654"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "2.130.76"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "7.607"
$ws.Range("E24").Value = "  -4.95%  "
$ws.Range("D25").Value = "161.48"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "8.975"
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("D27").Value = "0.1462"
$ws.Range("E27").Value = "  -5.84%  "
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").Value = "1.920"
$ws.Range("D30").Value = "1.375"
$ws.Range("E30").Value = "  -7.44%  "
$ws.Range("D31").Value = "4.416"
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").Value = "1.485"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").Value = "4.025"
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("D34").Value = "0.05187"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("D36").Value = "0.7088"
$ws.Range("E36").Value = "  -5.53%  "
$ws.Range("D37").Value = "0.9973"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").Value = "2.657"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "0.01858"
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").Value = "0.9133"
$ws.Range("E41").Value = "  +6.15%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4268"
$ws.Range("E42").Value = "  -6.26%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.895"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("D44").Value = "1.051.93"
$ws.Range("E44").Value = "  -6.80%  "
$ws.Range("D45").Value = "69.45"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "102.51"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.739"
$ws.Range("E48").Value = "  -6.65%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.127"
$ws.Range("E49").Value = "  -6.51%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.994.56"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.198"
$ws.Range("E51").Value = "  -3.57%  "
